$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price (D) and volume-change (E) values, and a few
# coin re-ordering updates (rows 38/39 and 50/51), matching the
# refreshed cryptos list published by the scheduled GitHub Action.

$ws.Range('D2').Value = '51.245.79'
$ws.Range('E2').Value = '  -15.56%  '
$ws.Range('D3').Value = '2.260.87'
$ws.Range('E3').Value = '  -22.13%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '431.75'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -18.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '115.89'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -19.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.451'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -17.84%  '
$ws.Range('D9').Value = '2.258.18'
$ws.Range('E9').Value = '  -22.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.13'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -14.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0829'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -22.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.295'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -17.86%  '
$ws.Range('E13').Value = '  -6.94%  '
$ws.Range('D14').Value = '2.643.47'
$ws.Range('E14').Value = '  -22.49%  '
$ws.Range('D15').Value = '51.483.15'
$ws.Range('E15').Value = '  -15.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.20'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -19.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000112'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -20.45%  '
$ws.Range('D18').Value = '2.263.48'
$ws.Range('E18').Value = '  -22.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.84'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -22.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '290.59'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -17.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.998'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.69'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.89%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.46'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -26.96%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.94'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -24.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.989'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '52.54'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -18.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.358'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -20.59%  '
$ws.Range('D28').Value = '2.340.65'
$ws.Range('E28').Value = '  -22.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.135'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -23.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.998'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.63'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -15.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '142.04'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -6.66%  '
$ws.Range('D33').Value = '0.0₃0614'
$ws.Range('E33').Value = '  -28.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '16.29'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -16.66%  '
$ws.Range('E35').Value = '  -23.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.55'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -18.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.993'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.41%  '
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.768'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -22.93%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.26'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -24.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.955'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -20.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '31.49'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -16.25%  '
$ws.Range('E42').Value = '  -2.17%  '
$ws.Range('E43').Value = '  -15.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0486'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -16.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.01'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -18.85%  '
$ws.Range('D46').Value = '1.853.26'
$ws.Range('E46').Value = '  -19.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.11'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -24.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0793'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -13.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0197'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -16.96%  '
$ws.Range('B50').Value = 'ZEEBU'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.62'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -5.26%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.85'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -22.07%  '

Write-Output "Applied cryptos list refresh"
